# Regenerate column G ("K") values for rows 2-33 with newly computed strikeout
# counts (s_vals), replacing the previous Strike# based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 7
    4  = 4
    5  = 6
    6  = 6
    7  = 10
    8  = 5
    9  = 6
    10 = 4
    11 = 7
    12 = 2
    13 = 10
    14 = 2
    15 = 7
    16 = 5
    17 = 6
    18 = 4
    19 = 6
    20 = 6
    21 = 5
    22 = 7
    23 = 2
    24 = 8
    25 = 8
    26 = 6
    27 = 3
    28 = 8
    29 = 8
    30 = 3
    31 = 3
    32 = 5
    33 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
